$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) and "Volume(1h)" (column E) values per row,
# mirroring the scheduled GitHub Actions refresh of the cryptos list.
# $null entries mean that column is unchanged for that row.
$updates = @{
    2  = @{ D = "62.127.85";  E = "  +1.51%  " }
    3  = @{ D = "2.417.66";   E = "  +1.64%  " }
    4  = @{ D = $null;        E = "  -0.09%  " }
    5  = @{ D = "558.62";     E = "  +1.68%  " }
    6  = @{ D = "143.32";     E = "  +2.97%  " }
    7  = @{ D = $null;        E = "  +0.05%  " }
    8  = @{ D = $null;        E = "  +0.62%  " }
    9  = @{ D = "2.414.65";   E = "  +1.46%  " }
    11 = @{ D = $null;        E = "  -0.98%  " }
    12 = @{ D = $null;        E = "  +1.07%  " }
    13 = @{ D = $null;        E = "  +1.15%  " }
    14 = @{ D = "26.23";      E = "  +4.43%  " }
    15 = @{ D = $null;        E = "  +5.46%  " }
    16 = @{ D = $null;        E = "  +2.54%  " }
    17 = @{ D = "61.999.73";  E = "  +1.42%  " }
    18 = @{ D = "2.417.72";   E = "  +1.29%  " }
    19 = @{ D = "11.18";      E = "  +2.87%  " }
    20 = @{ D = "4.21";       E = "  +1.13%  " }
    21 = @{ D = "324.59";     E = "  +0.93%  " }
    22 = @{ D = "6.76";       E = "  +0.21%  " }
    23 = @{ D = $null;        E = "  +0.05%  " }
    24 = @{ D = "65.44";      E = "  +1.62%  " }
    25 = @{ D = $null;        E = "  +1.10%  " }
    26 = @{ D = "9.03";       E = "  +7.22%  " }
    27 = @{ D = "594.16";     E = "  +17.04%  " }
    28 = @{ D = $null;        E = "  +0.07%  " }
    29 = @{ D = $null;        E = "  +1.33%  " }
    30 = @{ D = "0.0₃0940";   E = "  +5.63%  " }
    31 = @{ D = "8.31";       E = "  +1.53%  " }
    32 = @{ D = $null;        E = "  +5.10%  " }
    33 = @{ D = $null;        E = "  -1.57%  " }
    34 = @{ D = "1.88";       E = "  +2.20%  " }
    35 = @{ D = $null;        E = "  +2.60%  " }
    36 = @{ D = "5.71";       E = "  +5.45%  " }
    37 = @{ D = $null;        E = "  +0.02%  " }
    38 = @{ D = "4.79";       E = "  +2.40%  " }
    39 = @{ D = $null;        E = "  +1.37%  " }
    40 = @{ D = $null;        E = "  +0.33%  " }
    41 = @{ D = "150.73";     E = "  +2.76%  " }
    42 = @{ D = $null;        E = "  -2.48%  " }
    43 = @{ D = $null;        E = "  +0.06%  " }
    44 = @{ D = "2.36";       E = "  +12.77%  " }
    45 = @{ D = "151.33";     E = "  +1.67%  " }
    46 = @{ D = $null;        E = "  +1.55%  " }
    47 = @{ D = $null;        E = "  +3.40%  " }
    48 = @{ D = "20.28";      E = $null }
    49 = @{ D = "0.592";      E = "  +2.55%  " }
    50 = @{ D = $null;        E = "  +1.40%  " }
    51 = @{ D = $null;        E = "  +2.03%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($null -ne $vals.D) {
        # Force the "Price" cell to remain plain text (matches the source
        # data, which stores prices like "558.62" / "62.127.85" as text,
        # not numbers) by marking it Text before assigning, then restore
        # the default "Normal" style so no stray formatting is left behind.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $vals.D
        $dCell.Style = "Normal"
    }

    if ($null -ne $vals.E) {
        $ws.Cells.Item($row, 5).Value = $vals.E
    }
}
